$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.657.16"
$ws.Range("E2").Value = "  +1.23%  "

$ws.Range("D3").Value = "3.426.87"
$ws.Range("E3").Value = "  +0.19%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "585.30"
$ws.Range("E5").Value = "  +0.60%  "

$ws.Range("D6").Value = "181.65"
$ws.Range("E6").Value = "  +5.00%  "

$ws.Range("E7").Value = "  +6.25%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "3.426.21"
$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("D10").Value = "0.132"
$ws.Range("E10").Value = "  +2.65%  "

$ws.Range("D11").Value = "6.98"
$ws.Range("E11").Value = "  +2.12%  "

$ws.Range("E12").Value = "  +1.62%  "

$ws.Range("D13").Value = "4.025.86"
$ws.Range("E13").Value = "  +0.33%  "

$ws.Range("E14").Value = "  +0.90%  "

$ws.Range("D15").Value = "29.36"
$ws.Range("E15").Value = "  -1.44%  "

$ws.Range("D16").Value = "66.604.68"
$ws.Range("E16").Value = "  +0.97%  "

$ws.Range("D17").Value = "0.0000173"
$ws.Range("E17").Value = "  +2.18%  "

$ws.Range("D18").Value = "3.416.30"
$ws.Range("E18").Value = "  -0.11%  "

$ws.Range("D19").Value = "5.92"
$ws.Range("E19").Value = "  +0.76%  "

$ws.Range("D20").Value = "13.86"
$ws.Range("E20").Value = "  +1.34%  "

$ws.Range("D21").Value = "369.67"
$ws.Range("E21").Value = "  +1.40%  "

$ws.Range("D22").Value = "7.62"
$ws.Range("E22").Value = "  -0.54%  "

$ws.Range("D23").Value = "73.31"
$ws.Range("E23").Value = "  +3.02%  "

$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("E25").Value = "  +2.24%  "

$ws.Range("E26").Value = "  +7.01%  "

$ws.Range("D27").Value = "9.86"
$ws.Range("E27").Value = "  +2.36%  "

$ws.Range("E28").Value = "  +2.02%  "

$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("D30").Value = "5.82"
$ws.Range("E30").Value = "  +1.03%  "

$ws.Range("E31").Value = "  +1.34%  "

$ws.Range("D32").Value = "23.39"
$ws.Range("E32").Value = "  -1.11%  "

$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("D34").Value = "7.07"
$ws.Range("E34").Value = "  +0.84%  "

$ws.Range("D35").Value = "1.27"
$ws.Range("E35").Value = "  -1.34%  "

$ws.Range("E36").Value = "  +1.22%  "

$ws.Range("D37").Value = "163.30"
$ws.Range("E37").Value = "  +2.23%  "

$ws.Range("D38").Value = "0.869"
$ws.Range("E38").Value = "  -0.70%  "

$ws.Range("D39").Value = "27.64"
$ws.Range("E39").Value = "  -3.91%  "

$ws.Range("D40").Value = "1.81"
$ws.Range("E40").Value = "  +2.56%  "

$ws.Range("D41").Value = "2.66"
$ws.Range("E41").Value = "  +5.15%  "

$ws.Range("D42").Value = "4.43"
$ws.Range("E42").Value = "  +1.38%  "

$ws.Range("D43").Value = "2.708.53"
$ws.Range("E43").Value = "  +0.81%  "

$ws.Range("E44").Value = "  -0.66%  "

$ws.Range("D45").Value = "0.0690"
$ws.Range("E45").Value = "  +1.52%  "

$ws.Range("D46").Value = "25.06"
$ws.Range("E46").Value = "  +4.65%  "

$ws.Range("D47").Value = "338.50"
$ws.Range("E47").Value = "  +11.20%  "

$ws.Range("D48").Value = "39.95"
$ws.Range("E48").Value = "  +0.21%  "

$ws.Range("E49").Value = "  -0.02%  "

# Row 50/51: Arweave and Stellar swap ranking positions with new values
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.105"
$ws.Range("E50").Value = "  +4.05%  "

$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").Value = "32.20"
$ws.Range("E51").Value = "  +6.67%  "
